$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values in column D are plain numeric-looking
# strings (e.g. "310.08"). The source data stores these as text, so
# each such cell is switched to Text format before the value is written
# to stop Excel auto-converting it to a number.

$ws.Range("D2").Value = "43.626.09"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.386.31"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.08"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.58"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  -5.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.40"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "2.756.53"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.68"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "2.384.16"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "43.612.52"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E21").Value = "  -5.30%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.46"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.02"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.90"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.87"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.68"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.39"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +5.56%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.67"
$ws.Range("E40").Value = "  +8.13%  "
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  +13.11%  "
$ws.Range("D45").Value = "2.036.93"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.80"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0291"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.13"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "57.87"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("E51").Value = "  -0.17%  "
